# Applies updated profit-calculation values for Leve rows across all Job sheets.
# Generated from the upstream commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

# ALC (sheet 1), row 123
$ws = $wb.Worksheets.Item(1)
$ws.Range("H123").Value = 20567.273
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 20567.273
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 20567.273
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -30367.273

# ALC (sheet 1), row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 3041098.2
$ws.Range("I132").Value = 3664535
$ws.Range("K132").Value = 10993605
$ws.Range("M132").Value = -10991075

# ALC (sheet 1), row 137
$ws = $wb.Worksheets.Item(1)
$ws.Range("H137").Value = 2109.7368
$ws.Range("I137").Value = 1949.125
$ws.Range("J137").Value = 2966.3333
$ws.Range("K137").Value = 5847.375
$ws.Range("L137").Value = 8898.999899999999
$ws.Range("M137").Value = -3297.375
$ws.Range("N137").Value = -13998.9999

# ARM (sheet 2), row 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1293.1875
$ws.Range("I2").Value = 754.8889
$ws.Range("J2").Value = 1985.2858
$ws.Range("K2").Value = 754.8889
$ws.Range("L2").Value = 1985.2858
$ws.Range("M2").Value = -641.8889
$ws.Range("N2").Value = -2211.2858

# ARM (sheet 2), row 13
$ws = $wb.Worksheets.Item(2)
$ws.Range("H13").Value = 4000
$ws.Range("I13").Value = 7000
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = -6856
$ws.Range("N13").Value = -1288

# ARM (sheet 2), row 61
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 2343.8572
$ws.Range("J61").Value = 2801.4
$ws.Range("L61").Value = 2801.4
$ws.Range("N61").Value = -3225.4

# ARM (sheet 2), row 74
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 2375.9167
$ws.Range("I74").Value = 2062.9524
$ws.Range("J74").Value = 4566.6665
$ws.Range("K74").Value = 2062.9524
$ws.Range("L74").Value = 4566.6665
$ws.Range("M74").Value = -1188.9524
$ws.Range("N74").Value = -6314.6665

# ARM (sheet 2), row 76
$ws = $wb.Worksheets.Item(2)
$ws.Range("H76").Value = 27580
$ws.Range("J76").Value = 27580
$ws.Range("L76").Value = 27580
$ws.Range("N76").Value = -28256

# ARM (sheet 2), row 77
$ws = $wb.Worksheets.Item(2)
$ws.Range("H77").Value = 2375.9167
$ws.Range("I77").Value = 2062.9524
$ws.Range("J77").Value = 4566.6665
$ws.Range("K77").Value = 10314.762
$ws.Range("L77").Value = 22833.3325
$ws.Range("M77").Value = -5946.762000000001
$ws.Range("N77").Value = -31569.3325

# ARM (sheet 2), row 79
$ws = $wb.Worksheets.Item(2)
$ws.Range("H79").Value = 27580
$ws.Range("J79").Value = 27580
$ws.Range("L79").Value = 27580
$ws.Range("N79").Value = -29920

# ARM (sheet 2), row 116
$ws = $wb.Worksheets.Item(2)
$ws.Range("H116").Value = 1293.1875
$ws.Range("I116").Value = 754.8889
$ws.Range("J116").Value = 1985.2858
$ws.Range("K116").Value = 754.8889
$ws.Range("L116").Value = 1985.2858
$ws.Range("M116").Value = 1539.1111
$ws.Range("N116").Value = -6573.2858

# ARM (sheet 2), row 122
$ws = $wb.Worksheets.Item(2)
$ws.Range("H122").Value = 2278.1
$ws.Range("I122").Value = 1358.8572
$ws.Range("J122").Value = 2773.077
$ws.Range("K122").Value = 4076.5716
$ws.Range("L122").Value = 8319.231
$ws.Range("M122").Value = -1626.5716
$ws.Range("N122").Value = -13219.231

# ARM (sheet 2), row 132
$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value = 5175.9653
$ws.Range("I132").Value = 5351.174
$ws.Range("J132").Value = 4504.3335
$ws.Range("K132").Value = 16053.522
$ws.Range("L132").Value = 13513.0005
$ws.Range("M132").Value = -13523.522
$ws.Range("N132").Value = -18573.0005

# ARM (sheet 2), row 136
$ws = $wb.Worksheets.Item(2)
$ws.Range("H136").Value = 2343.8572
$ws.Range("J136").Value = 2801.4
$ws.Range("L136").Value = 8404.200000000001
$ws.Range("N136").Value = -13504.2

# BSM (sheet 3), row 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1293.1875
$ws.Range("I3").Value = 754.8889
$ws.Range("J3").Value = 1985.2858
$ws.Range("K3").Value = 754.8889
$ws.Range("L3").Value = 1985.2858
$ws.Range("M3").Value = -640.8889
$ws.Range("N3").Value = -2213.2858

# BSM (sheet 3), row 99
$ws = $wb.Worksheets.Item(3)
$ws.Range("H99").Value = 1566.6666
$ws.Range("I99").Value = 1566.6666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1566.6666
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -68.66660000000002
$ws.Range("N99").ClearContents()

# BSM (sheet 3), row 134
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 39116.965
$ws.Range("I134").Value = 68270.53
$ws.Range("J134").Value = 2675
$ws.Range("K134").Value = 204811.59
$ws.Range("L134").Value = 8025
$ws.Range("M134").Value = -202276.59
$ws.Range("N134").Value = -13095

# CRP (sheet 4), row 31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 5408112.5
$ws.Range("I31").Value = 2782.361
$ws.Range("J31").Value = 200000000
$ws.Range("K31").Value = 2782.361
$ws.Range("L31").Value = 200000000
$ws.Range("M31").Value = -2487.361
$ws.Range("N31").Value = -200000590

# CRP (sheet 4), row 34
$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 5408112.5
$ws.Range("I34").Value = 2782.361
$ws.Range("J34").Value = 200000000
$ws.Range("K34").Value = 2782.361
$ws.Range("L34").Value = 200000000
$ws.Range("M34").Value = -2580.361
$ws.Range("N34").Value = -200000404

# CRP (sheet 4), row 58
$ws = $wb.Worksheets.Item(4)
$ws.Range("H58").Value = 10417678
$ws.Range("I58").Value = 935.7083
$ws.Range("J58").Value = 41667904
$ws.Range("K58").Value = 935.7083
$ws.Range("L58").Value = 41667904
$ws.Range("M58").Value = -732.7083
$ws.Range("N58").Value = -41668310

# CRP (sheet 4), row 132
$ws = $wb.Worksheets.Item(4)
$ws.Range("H132").Value = 2145.973
$ws.Range("I132").Value = 1500.1333
$ws.Range("K132").Value = 4500.3999
$ws.Range("M132").Value = -1970.3999

# CRP (sheet 4), row 134
$ws = $wb.Worksheets.Item(4)
$ws.Range("H134").Value = 1278.75
$ws.Range("I134").Value = 1226.0667
$ws.Range("J134").Value = 1339.5385
$ws.Range("K134").Value = 3678.2001
$ws.Range("L134").Value = 4018.6155
$ws.Range("M134").Value = -1143.2001
$ws.Range("N134").Value = -9088.6155

# CRP (sheet 4), row 136
$ws = $wb.Worksheets.Item(4)
$ws.Range("H136").Value = 10417678
$ws.Range("I136").Value = 935.7083
$ws.Range("J136").Value = 41667904
$ws.Range("K136").Value = 2807.1249
$ws.Range("L136").Value = 125003712
$ws.Range("M136").Value = -257.1248999999998
$ws.Range("N136").Value = -125008812

# CUL (sheet 5), row 107
$ws = $wb.Worksheets.Item(5)
$ws.Range("H107").Value = 441
$ws.Range("I107").Value = 463.33334
$ws.Range("J107").Value = 427.6
$ws.Range("K107").Value = 1390.00002
$ws.Range("L107").Value = 1282.8
$ws.Range("M107").Value = 529.9999800000001
$ws.Range("N107").Value = -5122.8

# CUL (sheet 5), row 113
$ws = $wb.Worksheets.Item(5)
$ws.Range("H113").Value = 524.4792
$ws.Range("I113").Value = 525
$ws.Range("K113").Value = 1575
$ws.Range("M113").Value = 595

# CUL (sheet 5), row 132
$ws = $wb.Worksheets.Item(5)
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# CUL (sheet 5), row 136
$ws = $wb.Worksheets.Item(5)
$ws.Range("H136").Value = 1994.4897
$ws.Range("I136").Value = 1665
$ws.Range("J136").Value = 2008.5106
$ws.Range("K136").Value = 4995
$ws.Range("L136").Value = 6025.531800000001
$ws.Range("M136").Value = 105
$ws.Range("N136").Value = -16225.5318

# GSM (sheet 6), row 132
$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 107597.79
$ws.Range("I132").Value = 168555
$ws.Range("J132").Value = 3099.7144
$ws.Range("K132").Value = 505665
$ws.Range("L132").Value = 9299.143199999999
$ws.Range("M132").Value = -503135
$ws.Range("N132").Value = -14359.1432

# LTW (sheet 7), row 55
$ws = $wb.Worksheets.Item(7)
$ws.Range("H55").Value = 355.63635
$ws.Range("I55").Value = 350.2857
$ws.Range("J55").Value = 365
$ws.Range("K55").Value = 350.2857
$ws.Range("L55").Value = 365
$ws.Range("M55").Value = -177.2857
$ws.Range("N55").Value = -711

# LTW (sheet 7), row 132
$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value = 7015.528
$ws.Range("I132").Value = 8956.708000000001
$ws.Range("J132").Value = 3133.1667
$ws.Range("K132").Value = 26870.124
$ws.Range("L132").Value = 9399.500100000001
$ws.Range("M132").Value = -24340.124
$ws.Range("N132").Value = -14459.5001

# LTW (sheet 7), row 136
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 3654.6667
$ws.Range("I136").Value = 3777.8948
$ws.Range("J136").Value = 2985.7144
$ws.Range("K136").Value = 11333.6844
$ws.Range("L136").Value = 8957.143199999999
$ws.Range("M136").Value = -8783.6844
$ws.Range("N136").Value = -14057.1432

# WVR (sheet 8), row 132
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 1719.5
$ws.Range("I132").Value = 1078.8235
$ws.Range("J132").Value = 3897.8
$ws.Range("K132").Value = 3236.4705
$ws.Range("L132").Value = 11693.4
$ws.Range("M132").Value = -706.4704999999999
$ws.Range("N132").Value = -16753.4

# WVR (sheet 8), row 136
$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 5563.857
$ws.Range("J136").Value = 1385.3846
$ws.Range("L136").Value = 4156.1538
$ws.Range("N136").Value = -9256.1538
